$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0.02981232170111252
$ws.Range("B2").Value = 0.01030980947691171
$ws.Range("C2").Value = 0.07886553555727005
$ws.Range("D2").Value = 0.09581993626485329
$ws.Range("E2").Value = 0.1019920290467978
$ws.Range("F2").Value = 0.09880969240045873
$ws.Range("G2").Value = 0.0264176699616247
$ws.Range("H2").Value = "HIST"

# Row 3
$ws.Range("A3").Value = 0.02790758206592524
$ws.Range("B3").Value = 0.006167315309216211
$ws.Range("C3").Value = 0.002642247483078762
$ws.Range("D3").Value = 0.07064315185672221
$ws.Range("E3").Value = 0.1154623178720089
$ws.Range("F3").Value = 0.08765590895367681
$ws.Range("G3").Value = 0.01481525358142623

# Row 4
$ws.Range("A4").Value = 0.0578885796578331
$ws.Range("B4").Value = 0.02143430655904847
$ws.Range("C4").Value = 0.07993114739656448
$ws.Range("D4").Value = 0.1241097434545708
$ws.Range("E4").Value = 0.1286966308993251
$ws.Range("F4").Value = 0.1263615752190126
$ws.Range("G4").Value = 0.02338304711220633

# Row 5
$ws.Range("A5").Value = 0.00258757531719792
$ws.Range("B5").Value = 0.0008886521614781591
$ws.Range("C5").Value = 0.0002865124396746043
$ws.Range("D5").Value = 0.07399011370985233
$ws.Range("E5").Value = 0.0781784438921367
$ws.Range("F5").Value = 0.07602663841196711
$ws.Range("G5").Value = 0.0241775575112623

# Row 6
$ws.Range("A6").Value = 0.1868782761963731
$ws.Range("B6").Value = 0.07645071511630383
$ws.Range("C6").Value = 0.004124026745557785
$ws.Range("D6").Value = 0.2248082921625947
$ws.Range("E6").Value = 0.261170731249527
$ws.Range("F6").Value = 0.2416291371706875
$ws.Range("G6").Value = 0.01727697423717817

# New Row 7
$ws.Range("A7").Value = 0.2749036812750683
$ws.Range("B7").Value = 0.1212951601526189
$ws.Range("C7").Value = 0.02451329305768013
$ws.Range("D7").Value = 0.3094482230614211
$ws.Range("E7").Value = 0.3451583395644867
$ws.Range("F7").Value = 0.3263292516488247
$ws.Range("G7").Value = 0.01308993932176823
$ws.Range("H7").Value = "RESNET"
$ws.Range("I7").Value = "kmeans"
